$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '34.523.55'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.811.52'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.15%  '
Set-TextValue $ws.Range("D5") '228.63'
Set-TextValue $ws.Range("D6") '0.580'
$ws.Range("E6").Value = '  +4.38%  '
$ws.Range("E7").Value = '  +0.16%  '
Set-TextValue $ws.Range("D8") '35.86'
$ws.Range("E8").Value = '  +9.10%  '
$ws.Range("E9").Value = '  +2.41%  '
Set-TextValue $ws.Range("D10") '0.0700'
$ws.Range("E10").Value = '  +0.76%  '
Set-TextValue $ws.Range("D11") '0.0961'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '2.073.48'
$ws.Range("E12").Value = '  +0.81%  '
Set-TextValue $ws.Range("D13") '11.35'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").Value = '1.819.14'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("E15").Value = '  +1.70%  '
Set-TextValue $ws.Range("D16") '4.55'
$ws.Range("E16").Value = '  +5.37%  '
$ws.Range("D17").Value = '34.520.76'
$ws.Range("E17").Value = '  -0.02%  '
Set-TextValue $ws.Range("D18") '69.56'
$ws.Range("E18").Value = '  +0.53%  '
Set-TextValue $ws.Range("D19") '248.19'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").Value = '  -0.57%  '
Set-TextValue $ws.Range("D21") '11.63'
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("E22").Value = '  +0.26%  '
Set-TextValue $ws.Range("D23") '4.22'
$ws.Range("E23").Value = '  +1.04%  '
Set-TextValue $ws.Range("D24") '172.63'
$ws.Range("E24").Value = '  +1.03%  '
Set-TextValue $ws.Range("D25") '2.13'
$ws.Range("E25").Value = '  +2.52%  '
Set-TextValue $ws.Range("D26") '8.02'
$ws.Range("E26").Value = '  +8.71%  '
Set-TextValue $ws.Range("D27") '16.89'
$ws.Range("E27").Value = '  +1.41%  '
Set-TextValue $ws.Range("D28") '0.120'
$ws.Range("E28").Value = '  +3.89%  '
$ws.Range("E29").Value = '  -0.04%  '
Set-TextValue $ws.Range("D30") '4.10'
$ws.Range("E30").Value = '  +1.98%  '
Set-TextValue $ws.Range("D31") '3.89'
$ws.Range("E31").Value = '  +1.84%  '
Set-TextValue $ws.Range("D32") '0.0535'
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D35").Value = '1.406.33'
$ws.Range("E35").Value = '  -1.00%  '
Set-TextValue $ws.Range("D36") '0.683'
$ws.Range("E36").Value = '  -0.40%  '
Set-TextValue $ws.Range("D37") '2.52'
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("E39").Value = '  +0.73%  '
Set-TextValue $ws.Range("D40") '84.12'
$ws.Range("E40").Value = '  -0.43%  '
Set-TextValue $ws.Range("D41") '0.970'
$ws.Range("E41").Value = '  +2.20%  '
Set-TextValue $ws.Range("D42") '2.82'
$ws.Range("E42").Value = '  +1.28%  '
Set-TextValue $ws.Range("D43") '2.40'
$ws.Range("E43").Value = '  +0.54%  '
Set-TextValue $ws.Range("D44") '1.16'
$ws.Range("E44").Value = '  +5.93%  '
Set-TextValue $ws.Range("D45") '13.43'
$ws.Range("E45").Value = '  -4.20%  '
$ws.Range("D48").Value = '1.972.78'
$ws.Range("E48").Value = '  +0.78%  '
Set-TextValue $ws.Range("D49") '105.97'
$ws.Range("E49").Value = '  +0.52%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D46") '6.08'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D47") '0.0505'
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D50") '1.00'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0130'
$ws.Range("E51").Value = '  +0.58%  '
